$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F3").Value = 5544
$ws.Range("F9").Value = 7
$ws.Range("F12").Value = 1544
$ws.Range("F13").Value = 5089
$ws.Range("F15").Value = 240
$ws.Range("F16").Value = 214
$ws.Range("F17").Value = 28
$ws.Range("F18").Value = 10
$ws.Range("F19").Value = 111
$ws.Range("F20").Value = 4379
$ws.Range("F21").Value = 210
$ws.Range("F22").Value = 1159
$ws.Range("F24").Value = 66
$ws.Range("F25").Value = 211
$ws.Range("F27").Value = 173
$ws.Range("F32").Value = 13
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F4").Value = 5544
$ws.Range("F10").Value = 7
$ws.Range("F13").Value = 1544
$ws.Range("F14").Value = 5089
$ws.Range("F16").Value = 240
$ws.Range("F17").Value = 214
$ws.Range("F18").Value = 28
$ws.Range("F19").Value = 10
$ws.Range("F20").Value = 111
$ws.Range("F21").Value = 4379
$ws.Range("F22").Value = 210
$ws.Range("F23").Value = 1159
$ws.Range("F25").Value = 66
$ws.Range("F26").Value = 211
$ws.Range("F28").Value = 173
$ws.Range("F33").Value = 13
